{"js": "// \"make report look clean\"\n//\n// 1) The \"Solution\" paragraph had its sentence split across two runs\n//    (artifact of the old `_GoBack` edit-position bookmark sitting between\n//    them). Re-write that stretch of text in place so it collapses back\n//    into a single run, then drop the stale `_GoBack` bookmark.\n// 2) Add a blank spacer paragraph after the \"Motivation\" body paragraph\n//    (matching the blank paragraph style already used elsewhere in the\n//    section) so the sections are visually separated \u2014 the actual\n//    \"clean up\" of the report layout.\n// 3) The \"Related Cryptographic Concepts\" paragraph gets a fresh edit\n//    cursor parked right before \"block ciphers\u2026\", which is where Word's\n//    `_GoBack` bookmark now belongs (it always tracks the single most\n//    recent edit location in the document).\n\nconst body = context.document.body;\n\n// --- 1) Merge \"...which attem\" + \"pts to ... decode them\" back into one run ---\nconst solutionMatches = body.search(\n  \"which attempts to automatically find common encodings and tries to decode them\",\n  { matchCase: false }\n);\nsolutionMatches.load(\"items,text\");\nawait context.sync();\n\nif (solutionMatches.items.length > 0) {\n  const solutionRange = solutionMatches.items[0];\n  // Re-insert the identical text: this collapses the two runs (that used\n  // to be split by the old _GoBack bookmark) into a single run while\n  // leaving the trailing \".\" run (outside this range) untouched.\n  solutionRange.insertText(solutionRange.text, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) Insert a blank paragraph after the \"Motivation\" body paragraph ---\nconst motivationMatches = body.search(\n  \"Mini-CyberChef was developed with the goal of providing easier and more enjoyable access to cryptography by amateurs.\",\n  { matchCase: false }\n);\nmotivationMatches.load(\"items\");\nawait context.sync();\n\nif (motivationMatches.items.length > 0) {\n  const motivationParagraph = motivationMatches.items[0].paragraphs.getFirst();\n  motivationParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// --- 3) Split the \"Related Cryptographic Concepts\" sentence right before\n//        \"block ciphers\" and move the `_GoBack` bookmark there ---\n// A document only ever keeps a single `_GoBack` bookmark (it marks the\n// most recent edit location), so drop the old one before planting the\n// new one.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst conceptsMatches = body.search(\n  \"block ciphers and secure communication protocols.\",\n  { matchCase: false }\n);\nconceptsMatches.load(\"items\");\nawait context.sync();\n\nif (conceptsMatches.items.length > 0) {\n  const splitPoint = conceptsMatches.items[0].getRange(Word.RangeLocation.start);\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"make report look clean\"\n#\n# 1) The \"Solution\" paragraph's sentence was split across two runs\n#    (an artifact of the old `_GoBack` edit-position bookmark sitting\n#    between them, which is deleted along with the text it splits).\n#    Re-type that stretch of text in place so it collapses back into a\n#    single run.\n# 2) Add a blank spacer paragraph after the \"Motivation\" body paragraph\n#    (matching the blank paragraph style already used elsewhere in the\n#    section) so the sections are visually separated.\n# 3) Park a fresh edit cursor right before \"block ciphers\u2026\" in the\n#    \"Related Cryptographic Concepts\" paragraph \u2014 that's where Word's\n#    `_GoBack` bookmark now belongs, since it always tracks the single\n#    most recent edit location in the document.\n\n$d = $word.ActiveDocument\n\n# --- 1) Merge \"...which attem\" + \"pts to ... decode them\" back into one run ---\n$solutionRange = $d.Content\n$solutionRange.Find.ClearFormatting()\n$foundSolution = $solutionRange.Find.Execute(\"which attempts to automatically find common encodings and tries to decode them\")\nif ($foundSolution) {\n    $solutionText = $solutionRange.Text\n    # Delete + re-insert (rather than a plain re-assignment, which Word\n    # treats as a no-op when the text is unchanged) so the run split is\n    # actually rebuilt as a single run. This also removes the old\n    # `_GoBack` bookmark, since it sat inside this exact range.\n    $solutionRange.Delete()\n    $solutionRange.InsertBefore($solutionText)\n}\n\n# --- 2) Insert a blank paragraph after the \"Motivation\" body paragraph ---\n$motivationRange = $d.Content\n$motivationRange.Find.ClearFormatting()\n$foundMotivation = $motivationRange.Find.Execute(\"Mini-CyberChef was developed with the goal of providing easier and more enjoyable access to cryptography by amateurs.\")\nif ($foundMotivation) {\n    $motivationRange.Expand(4) | Out-Null  # wdParagraph - grow to the full paragraph\n    $motivationRange.InsertParagraphAfter()\n}\n\n# --- 3) Split the \"Related Cryptographic Concepts\" sentence right before\n#        \"block ciphers\" and move the `_GoBack` bookmark there ---\n# A document only ever keeps a single `_GoBack` bookmark (it marks the\n# most recent edit location), so drop the old one before planting the\n# new one.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$conceptsRange = $d.Content\n$conceptsRange.Find.ClearFormatting()\n$foundConcepts = $conceptsRange.Find.Execute(\"block ciphers and secure communication protocols.\")\nif ($foundConcepts) {\n    $splitPoint = $conceptsRange.Duplicate\n    $splitPoint.Collapse(1)  # wdCollapseStart\n    $d.Bookmarks.Add(\"_GoBack\", $splitPoint)\n}\n"}
